{"js": "// Fill in the \"Name\" / \"Student ID\" / \"Class\" submission-worksheet table\n// and merge the \"Task 2: \" + \"Odd and Even\" runs (dropping the stray\n// _GoBack bookmark) into a single run, as described by the commit.\n\n// --- 1. Locate the Name / Student ID / Class table (2nd table in body) ---\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst infoTable = tables.items[1];\n\n// Row 0 = Name, Row 1 = Student ID, Row 2 = Class; column 1 holds the value.\nconst nameCell = infoTable.getCell(0, 1);\nconst idCell = infoTable.getCell(1, 1);\nconst classCell = infoTable.getCell(2, 1);\n\nnameCell.body.paragraphs.load(\"items\");\nidCell.body.paragraphs.load(\"items\");\nclassCell.body.paragraphs.load(\"items\");\nawait context.sync();\n\n// Each value cell starts with a single empty paragraph \u2014 fill it in place\n// instead of inserting a brand-new one.\nnameCell.body.paragraphs.items[0].insertText(\"Muhammad Iylia Bin Mohd Hutta\", Word.InsertLocation.replace);\nidCell.body.paragraphs.items[0].insertText(\"P7474841\", Word.InsertLocation.replace);\nclassCell.body.paragraphs.items[0].insertText(\"NSDDA1/CE/2220/4\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 2. Merge \"Task 2: \" + \"Odd and Even\" into a single run/paragraph ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text === \"Task 2: Odd and Even\") {\n    para.insertText(\"Task 2: Odd and Even\", Word.InsertLocation.replace);\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fill in the \"Name\" / \"Student ID\" / \"Class\" submission-worksheet table\n# and merge the \"Task 2: \" + \"Odd and Even\" runs (dropping the stray\n# _GoBack bookmark) into a single run, as described by the commit.\n\n$d = $word.ActiveDocument\n\n# --- 1. Locate the Name / Student ID / Class table (2nd table in body) ---\n$infoTable = $d.Tables.Item(2)\n\n# Row 1 = Name, Row 2 = Student ID, Row 3 = Class; column 2 holds the value.\n$infoTable.Cell(1, 2).Range.Text = \"Muhammad Iylia Bin Mohd Hutta\"\n$infoTable.Cell(2, 2).Range.Text = \"P7474841\"\n$infoTable.Cell(3, 2).Range.Text = \"NSDDA1/CE/2220/4\"\n\n# --- 2. Merge \"Task 2: \" + \"Odd and Even\" into a single run/paragraph ---\n# (Find/Replace naturally collapses the split runs and removes the\n# now-orphaned _GoBack bookmark sitting between them.)\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Task 2: Odd and Even\"\n$find.Replacement.Text = \"Task 2: Odd and Even\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\nWrite-Output \"done\"\n"}
